$wb = $excel.ActiveWorkbook

# --- Defs ---
$ws = $wb.Worksheets.Item('Defs')
$ws.Range('A2').Value = 'lgs8j8fo-06px'
$ws.Range('B2').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C2').Value = 'lgs8j8fo'
$ws.Range('A3').Value = 'lgs8j8fo-0oj7'
$ws.Range('B3').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C3').Value = 'lgs8j8g8'
$ws.Range('A4').Value = 'lgs8j8fo-03ay'
$ws.Range('B4').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C4').Value = 'lgs8j8fo'
$ws.Range('A5').Value = 'lgs8j8g8-03n9'
$ws.Range('B5').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C5').Value = 'lgs8j8g8'

# --- Point Defs ---
$ws = $wb.Worksheets.Item('Point Defs')
$ws.Range('A2').Value = 'lgs8j8fo-pjps'
$ws.Range('B2').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C2').Value = 'lgs8j8fo'
$ws.Range('E2').Value = '0m7w'
$ws.Range('F2').Value = '8esq'
$ws.Range('G2').Value = 'Select Test'
$ws.Range('H2').Value = '⛏️'
$ws.Range('I2').Value = 'For testing selects'
$ws.Range('J2').Value = 'SELECT'
$ws.Range('K2').Value = 'COUNTOFEACH'
$ws.Range('A3').Value = 'lgs8j8fo-27z6i'
$ws.Range('B3').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C3').Value = 'lgs8j8fo'
$ws.Range('D3').Value = $false
$ws.Range('E3').Value = '05a8'
$ws.Range('F3').Value = '1vb5'
$ws.Range('G3').Value = 'Free Item'
$ws.Range('H3').Value = '🆓'
$ws.Range('J3').Value = 'TEXT'
$ws.Range('A4').Value = 'lgs8j8fo-62i6'
$ws.Range('B4').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C4').Value = 'lgs8j8g8'
$ws.Range('F4').Value = '0pc6'
$ws.Range('G4').Value = 'Numeric Thing'
$ws.Range('H4').Value = '#️⃣'
$ws.Range('I4').Value = 'Set a description'
$ws.Range('J4').Value = 'NUM'
$ws.Range('K4').Value = 'AVERAGE'
$ws.Range('A5').Value = 'lgs8j8fo-nljl'
$ws.Range('B5').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C5').Value = 'lgs8j8g8'
$ws.Range('D5').Value = $true
$ws.Range('F5').Value = '0tb7'
$ws.Range('G5').Value = 'Boolean Thing'
$ws.Range('H5').Value = '👍'
$ws.Range('I5').Value = 'Orig desc'
$ws.Range('J5').Value = 'BOOL'
$ws.Range('A6').Value = 'lgs8j8g8-012r'
$ws.Range('B6').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C6').Value = 'lgs8j8g8'
$ws.Range('F6').Value = '0pc6'
$ws.Range('G6').Value = 'Test Relabel'
$ws.Range('H6').Value = '#️⃣'
$ws.Range('I6').Value = 'Set a description'
$ws.Range('J6').Value = 'NUM'
$ws.Range('K6').Value = 'AVERAGE'
$ws.Range('A7').Value = 'lgs8j8g8-npea'
$ws.Range('B7').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C7').Value = 'lgs8j8g8'
$ws.Range('D7').Value = $false
$ws.Range('E7').Value = 'ay7l'
$ws.Range('F7').Value = '0tb7'
$ws.Range('G7').Value = 'Boolean Thing'
$ws.Range('H7').Value = '👎'
$ws.Range('I7').Value = 'Orig desc'
$ws.Range('J7').Value = 'BOOL'
$ws.Range('K7').Value = 'COUNT'

# --- Entry Base ---
$ws = $wb.Worksheets.Item('Entry Base')
$ws.Range('A2').Value = 'lgs8j8fo-0lfw'
$ws.Range('B2').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C2').Value = 'lgs8j8g8'
$ws.Range('A3').Value = 'lgs8j8fo-s0ps'
$ws.Range('B3').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C3').Value = 'lgs8j8fo'
$ws.Range('F3').Value = 'lgs8j8g0-mpib'
$ws.Range('G3').Value = '2023-04-22T12:09:10'
$ws.Range('A4').Value = 'lgs8j8g8-y87o'
$ws.Range('B4').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C4').Value = 'lgs8j8g8'

# --- Entry Points ---
$ws = $wb.Worksheets.Item('Entry Points')
$ws.Range('A2').Value = 'lgs8j8fo-afsz'
$ws.Range('B2').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C2').Value = 'lgs8j8g8'
$ws.Range('A3').Value = 'lgs8j8fo-x1oi'
$ws.Range('B3').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C3').Value = 'lgs8j8g8'
$ws.Range('A4').Value = 'lgs8j8g8-0eh5'
$ws.Range('B4').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C4').Value = 'lgs8j8g8'
$ws.Range('A5').Value = 'lgs8j8g8-u189'
$ws.Range('B5').Value = '2023-04-22T17:09:10.644Z'
$ws.Range('C5').Value = 'lgs8j8g8'

# --- Tag Defs ---
$ws = $wb.Worksheets.Item('Tag Defs')
$ws.Range('G1').ClearContents()
$ws.Range('H1').ClearContents()
$ws.Range('A2').Value = 'lgs8j8g0-063q'
$ws.Range('B2').Value = '2023-04-22T17:09:10.656Z'
$ws.Range('C2').Value = 'lgs8j8g0'
$ws.Range('D2').Value = $false
$ws.Range('E2').Value = '0q9d'
$ws.Range('F2').Value = 'My Tag!'
$ws.Range('A3').Value = 'lgs8j8g0-r9pi'
$ws.Range('B3').Value = '2023-04-22T17:09:10.656Z'
$ws.Range('C3').Value = 'lgs8j8g8'
$ws.Range('D3').Value = $true
$ws.Range('E3').Value = 'vvct'
$ws.Range('F3').Value = 'Orig Tag Label'
$ws.Range('A4').Value = 'lgs8j8g0-83ol'
$ws.Range('B4').Value = '2023-04-22T17:09:10.656Z'
$ws.Range('C4').Value = 'lgs8j8g0'
$ws.Range('D4').Value = $false
$ws.Range('E4').Value = '0vvi'
$ws.Range('F4').Value = 'Select Option Test'
$ws.Range('A5').Value = 'lgs8j8g8-kgcz'
$ws.Range('B5').Value = '2023-04-22T17:09:10.656Z'
$ws.Range('C5').Value = 'lgs8j8g8'
$ws.Range('D5').Value = $false
$ws.Range('E5').Value = 'vvct'
$ws.Range('F5').Value = 'New Label'

# --- Tags ---
$ws = $wb.Worksheets.Item('Tags')
$ws.Range('A2').Value = 'lgs8j8g0-31g6'
$ws.Range('B2').Value = '2023-04-22T17:09:10.656Z'
$ws.Range('C2').Value = 'lgs8j8g8'
$ws.Range('D2').Value = $true
$ws.Range('E2').Value = 'ay7l'
$ws.Range('G2').Value = 'vvct'
$ws.Range('A3').Value = 'lgs8j8g0-uq0p'
$ws.Range('B3').Value = '2023-04-22T17:09:10.656Z'
$ws.Range('C3').Value = 'lgs8j8g0'
$ws.Range('D3').Value = $false
$ws.Range('E3').Value = '0m7w'
$ws.Range('F3').Value = '8esq'
$ws.Range('G3').Value = '0vvi'
